$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$names = @(
    'Name',
    'Le''Veon Bell - PIT',
    'Todd Gurley - LAR',
    'Ezekiel Elliott - DAL',
    'Saquon Barkley - NYG',
    'David Johnson - ARI',
    'Melvin Gordon - LAC',
    'Leonard Fournette - JAC',
    'Dalvin Cook - MIN',
    'Alvin Kamara - NO',
    'Kareem Hunt - KC',
    'Christian McCaffrey - CAR',
    'Jordan Howard - CHI',
    'Marshawn Lynch - OAK',
    'Devonta Freeman - ATL',
    'Derrick Henry - TEN',
    'Sony Michel - NE',
    'Jerick McKinnon - SF',
    'Alex Collins - BAL',
    'LeSean McCoy - BUF',
    'Royce Freeman - DEN',
    'Lamar Miller - HOU',
    'Jay Ajayi - PHI',
    'Jamaal Williams - GB',
    'Ronald Jones II - TB',
    'Joe Mixon - CIN',
    'Kerryon Johnson - DET',
    'Kenyan Drake - MIA',
    'Aaron Jones - GB',
    'Rashaad Penny - SEA',
    'Mark Ingram - NO',
    'Marlon Mack - IND',
    'Isaiah Crowell - NYJ',
    'Rex Burkhead - NE',
    'Nyheim Hines - IND',
    'Dion Lewis - TEN',
    'Tevin Coleman - ATL',
    'C.J. Anderson - CAR',
    'Tarik Cohen - CHI',
    'Chris Thompson - WAS',
    'Carlos Hyde - CLE',
    'Duke Johnson - CLE',
    'Nick Chubb - CLE',
    'Devontae Booker - DEN',
    'Bilal Powell - NYJ',
    'Ty Montgomery - GB',
    'D''Onta Foreman - HOU',
    'Giovani Bernard - CIN',
    'Corey Clement - PHI',
    'Theo Riddick - DET',
    'Doug Martin - OAK',
    'Chris Carson - SEA',
    'Latavius Murray - MIN',
    'LeGarrette Blount - DET',
    'James White - NE',
    'Chris Ivory - BUF',
    'Peyton Barber - TB',
    'Matt Breida - SF',
    'Kenneth Dixon - BAL',
    'Frank Gore - MIA',
    'Austin Ekeler - LAC',
    'Samaje Perine - WAS',
    'Jordan Wilkins - IND',
    'Spencer Ware - KC',
    'T.J. Yeldon - JAC',
    'Wayne Gallman - NYG',
    'Ameer Abdullah - DET',
    'Javorius Allen - BAL',
    'Alfred Blue - HOU',
    'Jonathan Stewart - NYG',
    'Darren Sproles - PHI',
    'Kalen Ballage - MIA',
    'Elijah McGuire - NYJ',
    'Corey Grant - JAC',
    'Chase Edmonds - ARI',
    'C.J. Prosise - SEA',
    'John Kelly - LAR',
    'Rod Smith - DAL',
    'James Conner - PIT',
    'DeAndre Washington - OAK',
    'Jalen Richard - OAK',
    'Mike Gillislee - NE',
    'Terrance West - NO',
    'Tavon Austin - DAL',
    'Joe Williams - SF',
    'Justin Jackson - LAC',
    'De''Angelo Henderson - DEN',
    'Travaris Cadet - BUF',
    'Jeremy Hill - NE',
    'Charles Sims - TB',
    'Robert Turbin - IND',
    'Orleans Darkwa - ',
    'Mike Davis - SEA',
    'Mark Walton - CIN',
    'Adrian Peterson - ',
    'Jacquizz Rodgers - TB',
    'Thomas Rawls - NYJ',
    'Kyle Juszczyk - SF',
    'Cameron Artis-Payne - CAR',
    'Eddie Lacy - ',
    'Malcolm Brown - LAR',
    'Tion Green - IND',
    'Bo Scarbrough - DAL',
    'Charcandrick West - KC',
    'Shane Vereen - NO',
    'Ito Smith - ATL',
    'Paul Perkins - ',
    'Boston Scott - NO',
    'Jonathan Williams - NO',
    'Damien Williams - KC',
    'Jamaal Charles - ',
    'Benny Cunningham - CHI',
    'Wendell Smallwood - PHI',
    'Akrum Wadley - TEN',
    'Robert Kelley - WAS',
    'J.D. McKissic - SEA',
    'T.J. Logan - ARI',
    'Phillip Lindsay - DEN',
    'Elijhaa Penny - ARI',
    'Christine Michael - IND',
    'Tommy Bohanon - JAC',
    'Alfred Morris - SF',
    'Jeremy McNichols - SF',
    'Branden Oliver - IND',
    'Jamize Olawale - DAL',
    'Josh Adams - PHI',
    'Tyler Ervin - HOU',
    'Kerwynn Williams - KC',
    'Andy Janovich - DEN',
    'Mike Tolbert - ',
    'Aaron Ripkowski - GB',
    'Marcus Murphy - BUF',
    'Patrick DiMarco - BUF',
    'D.J. Foster - ARI',
    'Trey Edmunds - NO',
    'Mack Brown - MIN',
    'Roc Thomas - MIN',
    'David Fluellen - TEN',
    'CJ Ham - MIN',
    'Zach Line - NO',
    'Raheem Mostert - SF',
    'Byron Marshall - WAS',
    'Fitzgerald Toussaint - PIT',
    'Anthony Sherman - KC',
    'Stevan Ridley - PIT',
    'James Develin - NE',
    'Brian Hill - CIN',
    'Taiwan Jones - BUF',
    'Derek Watt - LAC',
    'Kapri Bibbs - WAS',
    'Dwayne Washington - DET',
    'Zach Zenner - DET',
    'David Williams - DEN',
    'Jalston Fowler - ',
    'Dan Vitale - CLE',
    'Roosevelt Nix - PIT',
    'Jay Prosch - HOU',
    'Matt Jones - PHI',
    'Elijah Hood - CAR',
    'Senorise Perry - MIA',
    'Trenton Cannon - NYJ',
    'Michael Burton - CHI',
    'Brandon Wilds - JAC',
    'Brandon Bolden - NE',
    'Devante Mays - GB',
    'Justin Davis - LAR',
    'Alex Armah - CAR',
    'Jeremy Langford - MIA',
    'Terron Ward - ',
    'Donnel Pumphrey - PHI',
    'Josh Ferguson - IND',
    'Terrence Magee - ATL',
    'Taquan Mizzell - CHI',
    'Matthew Dayes - CLE',
    'Daniel Lasco - NO',
    'Keith Marshall - '
)
# Clear the whole sheet's used range first to remove stale data
$ws.Cells.Clear()

# Write the new names into column B (which already carries the custom width formatting);
# then delete column A (the old Rank column) so column B shifts into column A,
# bringing its width formatting along with it.
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $names[$i]
}

$ws.Columns.Item(1).Delete()

# Final column A width (target OOXML width 36.140625; closest value this engine's
# sixth-of-a-character rounding can reach is 36.1666..., achieved with this input)
$ws.Columns.Item(1).ColumnWidth = 35.3

# Reset selection to A1 (no special selection override)
$ws.Range("A1").Select()
